# =========================================================================
# Edit script: restructure Input3G / Output3G data to the new breakdown
# (10 RSCP bands, 5 Ec/Io bands) per commit "Fixes for single excel template"
# =========================================================================

$wb = $excel.ActiveWorkbook
$wsIn  = $wb.Worksheets.Item("Input3G")
$wsOut = $wb.Worksheets.Item("Output3G")

# -------------------------------------------------------------------------
# 1) Input3G sheet - clear old content then write the new layout
# -------------------------------------------------------------------------
$wsIn.Range("A1:C24").ClearContents()

# --- RSCP header + 10 bands (rows 1-11) ---
$wsIn.Range("A1").Value = "UMTS RSCP (dBm)"

$wsIn.Range("A2").Value = "<= -110"
$wsIn.Range("B2").Value = 50954
$wsIn.Range("C2").Value = 0.29986876254259298

$wsIn.Range("A3").Value = "-105 to -110"
$wsIn.Range("B3").Value = 21024
$wsIn.Range("C3").Value = 0.12372808540439301

$wsIn.Range("A4").Value = "-100 to -105"
$wsIn.Range("B4").Value = 24977
$wsIn.Range("C4").Value = 0.146991837383254

$wsIn.Range("A5").Value = "-95 to -100"
$wsIn.Range("B5").Value = 21874
$wsIn.Range("C5").Value = 0.12873041001406499

$wsIn.Range("A6").Value = "-90 to -95"
$wsIn.Range("B6").Value = 15455
$wsIn.Range("C6").Value = [double]"9.0954031579381006E-2"

$wsIn.Range("A7").Value = "-85 to -90"
$wsIn.Range("B7").Value = 14374
$wsIn.Range("C7").Value = [double]"8.4592251693433998E-2"

$wsIn.Range("A8").Value = "-80 to -85"
$wsIn.Range("B8").Value = 12801
$wsIn.Range("C8").Value = [double]"7.5335008621653599E-2"

$wsIn.Range("A9").Value = "-75 to -80"
$wsIn.Range("B9").Value = 6249
$wsIn.Range("C9").Value = [double]"3.6775913512750003E-2"

$wsIn.Range("A10").Value = "-70 to -75"
$wsIn.Range("B10").Value = 2144
$wsIn.Range("C10").Value = [double]"1.2617628191924401E-2"

$wsIn.Range("A11").Value = "-15 to -70"
$wsIn.Range("B11").Value = 69
$wsIn.Range("C11").Value = [double]"4.06071056549808E-4"

# --- Ec/Io header + 5 bands (rows 13-18) ---
$wsIn.Range("A13").Value = "UMTS Ec/Io (dB)"

$wsIn.Range("A14").Value = "-34 to -13"
$wsIn.Range("B14").Value = 44872
$wsIn.Range("C14").Value = 0.26407565868844901

$wsIn.Range("A15").Value = "-13 to -10"
$wsIn.Range("B15").Value = 23066
$wsIn.Range("C15").Value = 0.135745434643157

$wsIn.Range("A16").Value = "-10 to -7"
$wsIn.Range("B16").Value = 30029
$wsIn.Range("C16").Value = 0.176723300828031

$wsIn.Range("A17").Value = "-7 to -4"
$wsIn.Range("B17").Value = 50225
$wsIn.Range("C17").Value = 0.29557853355382702

$wsIn.Range("A18").Value = "-4 to 0"
$wsIn.Range("B18").Value = 21729
$wsIn.Range("C18").Value = 0.127877072286533

# --- Band header + UMTS 900 / 2100 (rows 20-22) ---
$wsIn.Range("A20").Value = "UMTS Band"

$wsIn.Range("A21").Value = "UMTS 900"
$wsIn.Range("B21").Value = 73917
$wsIn.Range("C21").Value = 0.43500803314481401

$wsIn.Range("A22").Value = "UMTS 2100"
$wsIn.Range("B22").Value = 96004
$wsIn.Range("C22").Value = 0.56499196685518505

# --- trailing "Vodafone:..." label (row 24) ---
$wsIn.Range("A24").Value = "Vodafone:input\VF BMT Car pk UMTS.FMT"

# -------------------------------------------------------------------------
# 2) Output3G sheet - update labels/formulas in place (formatting untouched)
# -------------------------------------------------------------------------

# --- Serving / UMTS Ec/Io (dB) block, rows 2-13 ---
$wsOut.Range("B2").Value = "Serving "
$wsOut.Range("D2").Value = "UMTS  Ec/Io (dB)"

$wsOut.Range("B3").Value = "<= -110"
$wsOut.Range("D3").Formula = "=Input3G!C2"

$wsOut.Range("B4").Value = "-105 to -110"
$wsOut.Range("D4").Formula = "=Input3G!C3"

$wsOut.Range("B5").Value = "-100 to -105"
$wsOut.Range("D5").Formula = "=Input3G!C4"

$wsOut.Range("B6").Value = "-95 to -100"
$wsOut.Range("D6").Formula = "=Input3G!C5"

$wsOut.Range("B7").Value = "-90 to -95"
$wsOut.Range("D7").Formula = "=Input3G!C6"

$wsOut.Range("B8").Value = "-85 to -90"
$wsOut.Range("D8").Formula = "=Input3G!C7"

$wsOut.Range("B9").Value = "-80 to -85"
$wsOut.Range("D9").Formula = "=Input3G!C8"

$wsOut.Range("B10").Value = "-75 to -80"
$wsOut.Range("D10").Formula = "=Input3G!C9"

$wsOut.Range("B11").Value = "-70 to -75"
$wsOut.Range("D11").Formula = "=Input3G!C10"

$wsOut.Range("B12").Value = "-15 to -70"
$wsOut.Range("D12").Formula = "=Input3G!C11"

$wsOut.Range("B13").Formula = "=Input3G!A24"

# --- UMTS RSCP (dBm) / Quality Level block, rows 15-21 ---
$wsOut.Range("B15").Value = "UMTS  RSCP (dBm)"
$wsOut.Range("D15").Value = "Quality Level"

$wsOut.Range("B16").Value = "<= -13"
$wsOut.Range("D16").Formula = "=Input3G!C14"

$wsOut.Range("B17").Value = "-13 to -10"
$wsOut.Range("D17").Formula = "=Input3G!C15"

$wsOut.Range("B18").Value = "-10 to -7"
$wsOut.Range("D18").Formula = "=Input3G!C16"

$wsOut.Range("B19").Value = "-7 to -4"
$wsOut.Range("D19").Formula = "=Input3G!C17"

$wsOut.Range("B20").Value = "-4 to 0"
$wsOut.Range("D20").Formula = "=Input3G!C18"

$wsOut.Range("B21").Formula = "=Input3G!A24"

# --- Band / Serving block, rows 23-26 ---
$wsOut.Range("B23").Value = "Band"
$wsOut.Range("D23").Value = "Serving "

$wsOut.Range("B24").Value = "UMTS 900"
$wsOut.Range("D24").Formula = "=Input3G!C21"

$wsOut.Range("B25").Value = "UMTS 2100"
$wsOut.Range("D25").Formula = "=Input3G!C22"

$wsOut.Range("B26").Formula = "=Input3G!A24"

# -------------------------------------------------------------------------
# 3) Selection / active cell bookkeeping to mirror the diff
# -------------------------------------------------------------------------
$wsIn.Range("A24").Select()
$wsOut.Select()
$wsOut.Range("F20").Select()
